$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the two year-based columns (and their matching ".letter" columns)
# to the new "Kiefer (n. entb.)" / "Kiefer (entb.)" group labels.
$ws.Range("B1").Value = "Kiefer (n. entb.)"
$ws.Range("C1").Value = "Kiefer (entb.)"
$ws.Range("J1").Value = "Kiefer (n. entb.).letter"
$ws.Range("K1").Value = "Kiefer (entb.).letter"
